# Apply updated indicator quantile results:
#  - Column AJ (ScriptLatestRunVersion) text updated with new Git Commit ID
#  - Column AH (pid) values updated from 11992 to 17548
# Both changes apply uniformly to data rows 2 through 80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCommit = "IndicatorQuantiles.R, Git Commit ID: 24c5634628309d80791a95cb6332cf2c12927180"

for ($row = 2; $row -le 80; $row++) {
    $ws.Range("AH$row").Value = 17548
    $ws.Range("AJ$row").Value = $newCommit
}
